# Add 2022-Q3 data: insert a new "2022-Q3" sheet between "总计" and "2022-Q2",
# and add a corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (Total) sheet: shift existing rows down and insert the new
#    2022-Q3 row at the top of the data (row 2), pushing 2022-Q2 -> row 3
#    and 2022-Q1 -> row 4.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Preserve the bold/bordered style used by the index column (A2) by copying
# it down onto the newly-created A4 cell before overwriting its value.
$total.Range("A2").Copy($total.Range("A4"))

# Row 4 <- old row 3 (2022-Q1 : 2, 0.1)
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q1"
$total.Cells.Item(4, 3).Value = 2
$total.Cells.Item(4, 4).Value = 0.1

# Row 3 <- old row 2 (2022-Q2 : 5, 0.29)
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 5
$total.Cells.Item(3, 4).Value = 0.29

# Row 2 <- new 2022-Q3 data (3, 0.23)
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.23

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right before the existing "2022-Q2"
#    sheet, so the tab order becomes: 总计, 2022-Q3, 2022-Q2, 2022-Q1.
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

# Sheet indices shift by one now that $q3 has been inserted.
$q2 = $wb.Worksheets.Item(3)
$q1 = $wb.Worksheets.Item(4)

# Bring over the header row (text + bold/bordered style) and the per-row
# index-column style (A2:A4 bold/bordered) from the 2022-Q2 sheet, which has
# the identical layout/styling used by every quarterly fund-holdings sheet.
$q2.Range("B1:H1").Copy($q3.Range("B1:H1"))
$q2.Range("A2:H4").Copy($q3.Range("A2:H4"))

# Match the "2022-Q2" sheet's page margins (differs from the default used by
# a brand-new worksheet).
$q3.PageSetup.LeftMargin = 0.75 * 72
$q3.PageSetup.RightMargin = 0.75 * 72
$q3.PageSetup.TopMargin = 1 * 72
$q3.PageSetup.BottomMargin = 1 * 72
$q3.PageSetup.HeaderMargin = 0.5 * 72
$q3.PageSetup.FooterMargin = 0.5 * 72

# ---------------------------------------------------------------------------
# 3) Fill in the 2022-Q3 fund-holdings data.
#    Numeric-looking text columns (fund code + the 4 percentage/size columns)
#    must stay stored as text, matching the source data's convention, so we
#    toggle the number format to Text only while those values are written.
#    (Two single-area ranges are used instead of one multi-area "B2:B4,D2:G4"
#    range since multi-area NumberFormat assignment isn't applied reliably.)
# ---------------------------------------------------------------------------
$textColsCode = $q3.Range("B2:B4")
$textColsPct = $q3.Range("D2:G4")
$textColsCode.NumberFormat = "@"
$textColsPct.NumberFormat = "@"

# Row 2: fund 001672
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "001672"
$q3.Cells.Item(2, 3).Value = "国寿安保智慧生活股票"
$q3.Cells.Item(2, 4).Value = "4.26"
$q3.Cells.Item(2, 5).Value = "86.24"
$q3.Cells.Item(2, 6).Value = "3.00"
$q3.Cells.Item(2, 7).Value = "0.1278"
$q3.Cells.Item(2, 8).Value = 6

# Row 3: fund 004818
$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "004818"
$q3.Cells.Item(3, 3).Value = "国寿安保目标策略灵活配置混合A"
$q3.Cells.Item(3, 4).Value = "2.70"
$q3.Cells.Item(3, 5).Value = "45.00"
$q3.Cells.Item(3, 6).Value = "2.40"
$q3.Cells.Item(3, 7).Value = "0.0648"
$q3.Cells.Item(3, 8).Value = 3

# Row 4: fund 004819
$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "004819"
$q3.Cells.Item(4, 3).Value = "国寿安保目标策略灵活配置混合C"
$q3.Cells.Item(4, 4).Value = "1.73"
$q3.Cells.Item(4, 5).Value = "45.00"
$q3.Cells.Item(4, 6).Value = "2.40"
$q3.Cells.Item(4, 7).Value = "0.0415"
$q3.Cells.Item(4, 8).Value = 3

# Restore the General display format now that the text values are written,
# so the new sheet's data cells match the unstyled/no-numfmt cells used on
# every other quarterly sheet.
$textColsCode.ClearFormats()
$textColsPct.ClearFormats()

# ---------------------------------------------------------------------------
# 4) Keep "2022-Q1" as the selected/active tab, same as before the edit.
# ---------------------------------------------------------------------------
$q1.Activate()
